$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 5 through 16 (rows 1-4 remain unchanged)
$rows = @(
    @("Folha analítica -> Planilha", "Nome", "Profissional na folha analítica e fora da planilha"),
    @("Folha analítica -> Planilha", "Salário Liquido", "Profissional na folha analítica e fora da planilha"),
    @("Folha analítica -> Planilha", "Cargo", "Profissional na folha analítica e fora da planilha"),
    @("Folha analítica -> Planilha", "Data_Admissão", "Profissional na folha analítica e fora da planilha"),
    @("Folha analítica -> Planilha", "Matrícula", "Profissional na folha analítica e fora da planilha"),
    @("Folha analítica -> Comprovante de pagamento", "CLAYDERMAN DE OLIVEIRA FELISBERTO", "Erro no CPF ou matrícula"),
    @("Folha analítica -> Comprovante de pagamento", "ENEIDA SANTOS DE OLIVEIRA", "Erro no CPF ou matrícula"),
    @("Folha analítica -> Comprovante de pagamento", "JHONATAN MESQUITA DE ARAUJO", "Erro no CPF ou matrícula"),
    @("Folha analítica -> Comprovante de pagamento", "JULIANA MICHELLE ANIBAL ALVES", "Erro no CPF ou matrícula"),
    @("Comprovante de pagamento -> Folha analítica", "NOME CPF", "Profissional na folha analítica e fora da planilha"),
    @("Comprovante de pagamento -> Folha analítica", "VALOR", "Profissional na folha analítica e fora da planilha"),
    @("Planilha -> GFIP - SEFIP", "HAROLDO GONCALVES DA SILVA", "Profissional não encontrado pelo PIS/PASEP")
)

$startRow = 5
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
}

$wb.Save()
